$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1) # Overview
$ws2 = $wb.Worksheets.Item(2) # zh-cn
$ws3 = $wb.Worksheets.Item(3) # de-de

# Column widths are stored internally as (pixels+5)/6 where pixels is an
# integer pixel count derived from ColumnWidth (characters) rounded to the
# nearest whole pixel (6 px per character - Calibri 11 "Maximum Digit
# Width"). Use pre-computed ColumnWidth inputs that land on the desired
# stored widths.
$cwFor40 = 39.16666666666667      # -> stored width 40
$cwForC  = 29.166666666666664     # -> stored width 30 (closest grid point to 29.9777047293527)

# ---- Overview sheet: handback status text + column widths ----
$ws1.Range("E2").Value2 = "Handed back: in sync with en-US"
$ws1.Range("F2").Value2 = "Handed back: in sync with en-US"
$ws1.Columns.Item(5).ColumnWidth = $cwForC
$ws1.Columns.Item(6).ColumnWidth = $cwForC

# ---- zh-cn sheet ----
$ws2.Columns.Item(3).ColumnWidth = $cwForC
$ws2.Columns.Item(9).ColumnWidth = $cwFor40
$ws2.Columns.Item(10).ColumnWidth = $cwFor40
$ws2.Range("J2").Value2 = "7a25559b-172d-473d-abed-ba2bce1ca9ce.535f2aaa5ab51a347d02b4811ea6d3228d3001d3.zh-cn.xlf"
$ws2.Range("K2").Value2 = "2016-09-02 03:08:28"
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20ed137b76100d2ae245a1405ffff8d41f7b6310/e2e/7a25559b-172d-473d-abed-ba2bce1ca9ce.md", "", "", "7a25559b-172d-473d-abed-ba2bce1ca9ce.md")

# ---- de-de sheet ----
$ws3.Columns.Item(3).ColumnWidth = $cwForC
$ws3.Columns.Item(9).ColumnWidth = $cwFor40
$ws3.Columns.Item(10).ColumnWidth = $cwFor40
$ws3.Range("J2").Value2 = "7a25559b-172d-473d-abed-ba2bce1ca9ce.535f2aaa5ab51a347d02b4811ea6d3228d3001d3.de-de.xlf"
$ws3.Range("K2").Value2 = "2016-09-02 03:08:35"
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20ed137b76100d2ae245a1405ffff8d41f7b6310/e2e/7a25559b-172d-473d-abed-ba2bce1ca9ce.md", "", "", "7a25559b-172d-473d-abed-ba2bce1ca9ce.md")
